# Update the "Förändrad" (Changed) date column (column C) for all data
# rows (rows 2 through 120) from 45184 (2023-09-15) to 45185 (2023-09-16).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 120; $r++) {
    $ws.Cells.Item($r, 3).Value = 45185
}
